$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

# Column A: date stored as literal text (matches existing rows' format),
# not auto-converted to a date serial number.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/24/2025"
$dateCell.Style = "Normal"

# Columns B and C: numeric allocation values.
$ws.Cells.Item($row, 2).Value = 0.1958495701456533
$ws.Cells.Item($row, 3).Value = 0.8041504298543467
